$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 42.2
$ws.Range("I8").Value = 42.2
$ws.Range("K8").Value = 126.6
$ws.Range("M8").Value = 12.39999999999999
$ws.Range("H98").Value = 33944.293
$ws.Range("J98").Value = 66822.75
$ws.Range("L98").Value = 66822.75
$ws.Range("N98").Value = -69818.75
$ws.Range("H122").Value = 33944.293
$ws.Range("J122").Value = 66822.75
$ws.Range("L122").Value = 200468.25
$ws.Range("N122").Value = -205368.25
$ws.Range("H136").Value = 39160.625
$ws.Range("J136").Value = 39160.625
$ws.Range("L136").Value = 39160.625
$ws.Range("N136").Value = -49360.625
$ws.Range("H137").Value = 2710.1309
$ws.Range("I137").Value = 741.7931
$ws.Range("K137").Value = 2225.3793
$ws.Range("M137").Value = 324.6206999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3803.77
$ws.Range("I32").Value = 3522.677
$ws.Range("J32").Value = 10550
$ws.Range("K32").Value = 3522.677
$ws.Range("L32").Value = 10550
$ws.Range("M32").Value = -3235.677
$ws.Range("N32").Value = -11124
$ws.Range("H61").Value = 3165.2903
$ws.Range("I61").Value = 2069.4443
$ws.Range("J61").Value = 4682.615
$ws.Range("K61").Value = 2069.4443
$ws.Range("L61").Value = 4682.615
$ws.Range("M61").Value = -1857.4443
$ws.Range("N61").Value = -5106.615
$ws.Range("H102").Value = 20089.182
$ws.Range("I102").Value = 1602
$ws.Range("J102").Value = 35495.168
$ws.Range("K102").Value = 1602
$ws.Range("L102").Value = 35495.168
$ws.Range("M102").Value = 20
$ws.Range("N102").Value = -38739.168
$ws.Range("H110").Value = 1488.5834
$ws.Range("I110").Value = 1085
$ws.Range("K110").Value = 1085
$ws.Range("M110").Value = 960
$ws.Range("H122").Value = 1632.2858
$ws.Range("I122").Value = 1665.5385
$ws.Range("K122").Value = 4996.6155
$ws.Range("M122").Value = -2546.6155
$ws.Range("H132").Value = 8476249
$ws.Range("I132").Value = 12196228
$ws.Range("J132").Value = 2964.2778
$ws.Range("K132").Value = 36588684
$ws.Range("L132").Value = 8892.8334
$ws.Range("M132").Value = -36586154
$ws.Range("N132").Value = -13952.8334
$ws.Range("H136").Value = 3165.2903
$ws.Range("I136").Value = 2069.4443
$ws.Range("J136").Value = 4682.615
$ws.Range("K136").Value = 6208.3329
$ws.Range("L136").Value = 14047.845
$ws.Range("M136").Value = -3658.3329
$ws.Range("N136").Value = -19147.845

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 558.4
$ws.Range("I94").Value = 454.85715
$ws.Range("J94").Value = 800
$ws.Range("K94").Value = 454.85715
$ws.Range("L94").Value = 800
$ws.Range("M94").Value = -3.85714999999999
$ws.Range("N94").Value = -1702
$ws.Range("H105").Value = 2495.2144
$ws.Range("I105").Value = 2546.4546
$ws.Range("K105").Value = 2546.4546
$ws.Range("M105").Value = -799.4546
$ws.Range("H107").Value = 1812.1482
$ws.Range("J107").Value = 2375.1667
$ws.Range("L107").Value = 2375.1667
$ws.Range("N107").Value = -6215.1667
$ws.Range("H134").Value = 3883.9473
$ws.Range("I134").Value = 3092.3333
$ws.Range("K134").Value = 9276.999899999999
$ws.Range("M134").Value = -6741.999899999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3829.141
$ws.Range("I31").Value = 1464.5
$ws.Range("J31").Value = 5256.849
$ws.Range("K31").Value = 1464.5
$ws.Range("L31").Value = 5256.849
$ws.Range("M31").Value = -1169.5
$ws.Range("N31").Value = -5846.849
$ws.Range("H34").Value = 3829.141
$ws.Range("I34").Value = 1464.5
$ws.Range("J34").Value = 5256.849
$ws.Range("K34").Value = 1464.5
$ws.Range("L34").Value = 5256.849
$ws.Range("M34").Value = -1262.5
$ws.Range("N34").Value = -5660.849
$ws.Range("H58").Value = 4388.2617
$ws.Range("J58").Value = 2647.8572
$ws.Range("L58").Value = 2647.8572
$ws.Range("N58").Value = -3053.8572
$ws.Range("H63").Value = 30000
$ws.Range("J63").Value = 30000
$ws.Range("L63").Value = 30000
$ws.Range("N63").Value = -31372
$ws.Range("H66").Value = 30000
$ws.Range("J66").Value = 30000
$ws.Range("L66").Value = 90000
$ws.Range("N66").Value = -96864
$ws.Range("H112").Value = 37222.332
$ws.Range("J112").Value = 37222.332
$ws.Range("L112").Value = 37222.332
$ws.Range("N112").Value = -40176.332
$ws.Range("H132").Value = 46726.188
$ws.Range("I132").Value = 1768.1666
$ws.Range("J132").Value = 104529.36
$ws.Range("K132").Value = 5304.4998
$ws.Range("L132").Value = 313588.08
$ws.Range("M132").Value = -2774.4998
$ws.Range("N132").Value = -318648.08
$ws.Range("H136").Value = 4388.2617
$ws.Range("J136").Value = 2647.8572
$ws.Range("L136").Value = 7943.571599999999
$ws.Range("N136").Value = -13043.5716

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H88").Value = 7010.6665
$ws.Range("J88").Value = 7010.6665
$ws.Range("L88").Value = 21031.9995
$ws.Range("N88").Value = -21887.9995
$ws.Range("H91").Value = 7010.6665
$ws.Range("J91").Value = 7010.6665
$ws.Range("L91").Value = 21031.9995
$ws.Range("N91").Value = -23995.9995

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 4819.1
$ws.Range("I97").Value = 1566.2727
$ws.Range("K97").Value = 1566.2727
$ws.Range("M97").Value = -1070.2727
$ws.Range("H102").Value = 3511.4167
$ws.Range("I102").Value = 3613.6667
$ws.Range("J102").Value = 3204.6667
$ws.Range("K102").Value = 3613.6667
$ws.Range("L102").Value = 3204.6667
$ws.Range("M102").Value = -1991.6667
$ws.Range("N102").Value = -6448.6667
$ws.Range("H113").Value = 7458.1665
$ws.Range("I113").Value = 10168.417
$ws.Range("K113").Value = 10168.417
$ws.Range("M113").Value = -7998.416999999999
$ws.Range("H126").Value = 2409.2173
$ws.Range("I126").Value = 3482.4
$ws.Range("J126").Value = 2111.111
$ws.Range("K126").Value = 10447.2
$ws.Range("L126").Value = 6333.333
$ws.Range("M126").Value = -7977.200000000001
$ws.Range("N126").Value = -11273.333

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 6099.0625
$ws.Range("I132").Value = 10119.4
$ws.Range("J132").Value = 4271.636
$ws.Range("K132").Value = 30358.2
$ws.Range("L132").Value = 12814.908
$ws.Range("M132").Value = -27828.2
$ws.Range("N132").Value = -17874.908

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H33").Value = 14293.667
$ws.Range("I33").Value = 8900
$ws.Range("J33").Value = 16990.5
$ws.Range("K33").Value = 8900
$ws.Range("L33").Value = 16990.5
$ws.Range("M33").Value = -8650
$ws.Range("N33").Value = -17490.5
$ws.Range("H36").Value = 14293.667
$ws.Range("I36").Value = 8900
$ws.Range("J36").Value = 16990.5
$ws.Range("K36").Value = 8900
$ws.Range("L36").Value = 16990.5
$ws.Range("M36").Value = -8650
$ws.Range("N36").Value = -17490.5
$ws.Range("H81").Value = 1772.3334
$ws.Range("I81").Value = 1742.8572
$ws.Range("J81").Value = 1791.091
$ws.Range("K81").Value = 3485.7144
$ws.Range("L81").Value = 3582.182
$ws.Range("M81").Value = -2424.7144
$ws.Range("N81").Value = -5704.182
$ws.Range("H84").Value = 1772.3334
$ws.Range("I84").Value = 1742.8572
$ws.Range("J84").Value = 1791.091
$ws.Range("K84").Value = 17428.572
$ws.Range("L84").Value = 17910.91
$ws.Range("M84").Value = -12124.572
$ws.Range("N84").Value = -28518.91
$ws.Range("H122").Value = 1786658
$ws.Range("I122").Value = 2198579.5
$ws.Range("K122").Value = 6595738.5
$ws.Range("M122").Value = -6593288.5
$ws.Range("H126").Value = 2451491.5
$ws.Range("I126").Value = 2674272.5
$ws.Range("J126").Value = 900
$ws.Range("K126").Value = 8022817.5
$ws.Range("L126").Value = 2700
$ws.Range("M126").Value = -8020347.5
$ws.Range("N126").Value = -7640
$ws.Range("H132").Value = 1834.1333
$ws.Range("I132").Value = 1228.5454
$ws.Range("J132").Value = 3499.5
$ws.Range("K132").Value = 3685.6362
$ws.Range("L132").Value = 10498.5
$ws.Range("M132").Value = -1155.6362
$ws.Range("N132").Value = -15558.5
$ws.Range("H136").Value = 17010.666
$ws.Range("I136").Value = 36338.098
$ws.Range("J136").Value = 2397.244
$ws.Range("K136").Value = 109014.294
$ws.Range("L136").Value = 7191.732
$ws.Range("M136").Value = -106464.294
$ws.Range("N136").Value = -12291.732
